$wb = $excel.ActiveWorkbook

# Hunk 1: sheet 'Option 1 - LR1 - DN1 (80-20)', row 22
$ws = $wb.Worksheets.Item('Option 1 - LR1 - DN1 (80-20)')
$ws.Range('B22').Value = 'Random'
$ws.Range('C22').Value = '{''module__num_units'': 10, ''module__activation_func'': ReLU(), ''lr'': 0.1}'
$ws.Range('D22').Value = 0.07838523387908936
$ws.Range('E22').Value = 0.215093806385994
$ws.Range('F22').Value = 2.577588558197021
$ws.Range('G22').Value = 0.2799736306852654
$ws.Range('H22').Value = 40.89487195014954
$ws.Range('J22').Value = 'Random'
$ws.Range('K22').Value = '{''module__num_units'': 10, ''module__activation_func'': ReLU(), ''lr'': 0.1}'
$ws.Range('L22').Value = 0.08215077221393585
$ws.Range('M22').Value = 0.2195072770118713
$ws.Range('N22').Value = 2.674199819564819
$ws.Range('O22').Value = 0.2866195600686315
$ws.Range('P22').Value = 40.85392653942108

# Hunk 2: sheet 'Option 1 - LR2 - DN2 (60-40)', row 22
$ws = $wb.Worksheets.Item('Option 1 - LR2 - DN2 (60-40)')
$ws.Range('B22').Value = 'Random'
$ws.Range('C22').Value = '{''module__num_units'': 10, ''module__activation_func'': ReLU(), ''lr'': 0.1}'
$ws.Range('D22').Value = 220.8598022460938
$ws.Range('E22').Value = 14.85894680023193
$ws.Range('F22').Value = 67.21656799316406
$ws.Range('G22').Value = 14.86135263850817
$ws.Range('H22').Value = 185.0137591362
$ws.Range('J22').Value = 'Random'
$ws.Range('K22').Value = '{''module__num_units'': 10, ''module__activation_func'': Tanh(), ''lr'': 0.02}'
$ws.Range('L22').Value = 238.7121276855469
$ws.Range('M22').Value = 15.44793796539307
$ws.Range('N22').Value = 69.78421783447266
$ws.Range('O22').Value = 15.45031157244238
$ws.Range('P22').Value = 185.5399012565613

# Hunk 3: sheet 'Option 1 - LR2 - DN2 (70-30)', row 22
$ws = $wb.Worksheets.Item('Option 1 - LR2 - DN2 (70-30)')
$ws.Range('B22').Value = 'Random'
$ws.Range('C22').Value = '{''module__num_units'': 20, ''module__activation_func'': ReLU(), ''lr'': 0.02}'
$ws.Range('D22').Value = 225.9319763183594
$ws.Range('E22').Value = 15.00244426727295
$ws.Range('F22').Value = 7.622591495513916
$ws.Range('G22').Value = 15.03103377410747
$ws.Range('H22').Value = 151.177179813385
$ws.Range('J22').Value = 'Random'
$ws.Range('K22').Value = '{''module__num_units'': 10, ''module__activation_func'': ReLU(), ''lr'': 0.02}'
$ws.Range('L22').Value = 133.8338623046875
$ws.Range('M22').Value = 11.53218555450439
$ws.Range('N22').Value = 5.922179222106934
$ws.Range('O22').Value = 11.56865862166775
$ws.Range('P22').Value = 140.9954428672791

# Hunk 4: sheet 'Option 1 - LR2 - DN2 (80-20)', row 22
$ws = $wb.Worksheets.Item('Option 1 - LR2 - DN2 (80-20)')
$ws.Range('B22').Value = 'Random'
$ws.Range('C22').Value = '{''module__num_units'': 50, ''module__activation_func'': Tanh(), ''lr'': 0.1}'
$ws.Range('D22').Value = 221.8839416503906
$ws.Range('E22').Value = 13.55358409881592
$ws.Range('F22').Value = 0.919352650642395
$ws.Range('G22').Value = 14.89576925339509
$ws.Range('H22').Value = 170.1185345649719
$ws.Range('J22').Value = 'Random'
$ws.Range('K22').Value = '{''module__num_units'': 50, ''module__activation_func'': Tanh(), ''lr'': 0.1}'
$ws.Range('L22').Value = 221.9650421142578
$ws.Range('M22').Value = 13.55808067321777
$ws.Range('N22').Value = 0.9219601154327393
$ws.Range('O22').Value = 14.89849126973123
$ws.Range('P22').Value = 170.1954245567322

# Hunk 5: sheet 'Option 1 - NLR1 - DN1 (70-30)', row 22
$ws = $wb.Worksheets.Item('Option 1 - NLR1 - DN1 (70-30)')
$ws.Range('B22').Value = 'Random'
$ws.Range('C22').Value = '{''module__num_units'': 50, ''module__activation_func'': Tanh(), ''lr'': 0.1}'
$ws.Range('D22').Value = 0.08294864743947983
$ws.Range('E22').Value = 0.2235938757658005
$ws.Range('F22').Value = 29775519285248
$ws.Range('G22').Value = 0.2880080683583011
$ws.Range('H22').Value = 43.07083785533905
$ws.Range('J22').Value = 'Random'
$ws.Range('K22').Value = '{''module__num_units'': 50, ''module__activation_func'': Tanh(), ''lr'': 0.1}'
$ws.Range('L22').Value = 0.08614324778318405
$ws.Range('M22').Value = 0.2296033203601837
$ws.Range('N22').Value = 30235582005248
$ws.Range('O22').Value = 0.2935016997960728
$ws.Range('P22').Value = 43.15841794013977

# Hunk 6: sheet 'Option 1 - NLR1 - DN2 (70-30)', row 22
$ws = $wb.Worksheets.Item('Option 1 - NLR1 - DN2 (70-30)')
$ws.Range('B22').Value = 'Random'
$ws.Range('C22').Value = '{''module__num_units'': 20, ''module__activation_func'': Tanh(), ''lr'': 0.1}'
$ws.Range('D22').Value = 1.209032416343689
$ws.Range('E22').Value = 0.5945903062820435
$ws.Range('F22').Value = 0.7156850695610046
$ws.Range('G22').Value = 1.099560101287642
$ws.Range('H22').Value = 26.51034593582153
$ws.Range('J22').Value = 'Random'
$ws.Range('K22').Value = '{''module__num_units'': 20, ''module__activation_func'': Tanh(), ''lr'': 0.1}'
$ws.Range('L22').Value = 1.169228196144104
$ws.Range('M22').Value = 0.6135646104812622
$ws.Range('N22').Value = 0.706624448299408
$ws.Range('O22').Value = 1.081308557324922
$ws.Range('P22').Value = 30.0974428653717

# Hunk 7: sheet 'Option 1 - NLR2 - DN1 (70-30)', row 22
$ws = $wb.Worksheets.Item('Option 1 - NLR2 - DN1 (70-30)')
$ws.Range('B22').Value = 'Random'
$ws.Range('C22').Value = '{''module__num_units'': 10, ''module__activation_func'': ReLU(), ''lr'': 0.1}'
$ws.Range('D22').Value = 0.08387590944766998
$ws.Range('E22').Value = 0.2306056022644043
$ws.Range('F22').Value = 2.468596935272217
$ws.Range('G22').Value = 0.2896133792621984
$ws.Range('H22').Value = 41.43242239952087

# Hunk 8: sheet 'Option 1 - NLR2 - DN2 (70-30)', row 22
$ws = $wb.Worksheets.Item('Option 1 - NLR2 - DN2 (70-30)')
$ws.Range('B22').Value = 'Random'
$ws.Range('C22').Value = '{''module__num_units'': 10, ''module__activation_func'': ReLU(), ''lr'': 0.1}'
$ws.Range('D22').Value = 0.9527064561843872
$ws.Range('E22').Value = 0.5113569498062134
$ws.Range('F22').Value = 0.5067695379257202
$ws.Range('G22').Value = 0.9760668297736519
$ws.Range('H22').Value = 25.04952549934387

# Hunk 9: sheet 'Option 1 - LR1 - DN1 (60-40)', row 22
$ws = $wb.Worksheets.Item('Option 1 - LR1 - DN1 (60-40)')
$ws.Range('B22').Value = 'Random'
$ws.Range('C22').Value = '{''module__num_units'': 50, ''module__activation_func'': Tanh(), ''lr'': 0.1}'
$ws.Range('D22').Value = 0.07530313730239868
$ws.Range('E22').Value = 0.2082056403160095
$ws.Range('F22').Value = 2.210495471954346
$ws.Range('G22').Value = 0.2744141711034594
$ws.Range('H22').Value = 41.23804569244385
$ws.Range('J22').Value = 'Random'
$ws.Range('K22').Value = '{''module__num_units'': 20, ''module__activation_func'': Tanh(), ''lr'': 0.1}'
$ws.Range('L22').Value = 0.07020536065101624
$ws.Range('M22').Value = 0.2018416970968246
$ws.Range('N22').Value = 2.061851263046265
$ws.Range('O22').Value = 0.2649629420334403
$ws.Range('P22').Value = 41.15466177463531

# Hunk 10: sheet 'Option 1 - LR1 - DN1 (70-30)', row 22
$ws = $wb.Worksheets.Item('Option 1 - LR1 - DN1 (70-30)')
$ws.Range('B22').Value = 'Random'
$ws.Range('C22').Value = '{''module__num_units'': 20, ''module__activation_func'': ReLU(), ''lr'': 0.1}'
$ws.Range('D22').Value = 0.07684499770402908
$ws.Range('E22').Value = 0.2159523963928223
$ws.Range('F22').Value = 2.456914663314819
$ws.Range('G22').Value = 0.2772093030618364
$ws.Range('H22').Value = 40.13096690177917
$ws.Range('J22').Value = 'Random'
$ws.Range('K22').Value = '{''module__num_units'': 20, ''module__activation_func'': Tanh(), ''lr'': 0.1}'
$ws.Range('L22').Value = 0.07276012003421783
$ws.Range('M22').Value = 0.2020772397518158
$ws.Range('N22').Value = 2.508364915847778
$ws.Range('O22').Value = 0.2697408386474281
$ws.Range('P22').Value = 40.30642807483673

# Hunk 11: sheet 'Option 1 - LR1 - DN2 (80-20)', row 23
$ws = $wb.Worksheets.Item('Option 1 - LR1 - DN2 (80-20)')
$ws.Range('B23').Value = 'Random'
$ws.Range('C23').Value = '{''module__num_units'': 20, ''module__activation_func'': Tanh(), ''lr'': 0.1}'
$ws.Range('D23').Value = 0.9624125957489014
$ws.Range('E23').Value = 0.5236716270446777
$ws.Range('F23').Value = 0.4887340068817139
$ws.Range('G23').Value = 0.9810262971750051
$ws.Range('H23').Value = 25.3871887922287
$ws.Range('J23').Value = 'Random'
$ws.Range('K23').Value = '{''module__num_units'': 50, ''module__activation_func'': Tanh(), ''lr'': 0.1}'
$ws.Range('L23').Value = 0.9571841955184937
$ws.Range('M23').Value = 0.5163381695747375
$ws.Range('N23').Value = 0.4868794977664948
$ws.Range('O23').Value = 0.9783579076792366
$ws.Range('P23').Value = 25.36321878433228

# Hunk 12: sheet 'Option 1 - LR1 - DN2 (60-40)', row 23
$ws = $wb.Worksheets.Item('Option 1 - LR1 - DN2 (60-40)')
$ws.Range('B23').Value = 'Random'
$ws.Range('C23').Value = '{''module__num_units'': 10, ''module__activation_func'': ReLU(), ''lr'': 0.1}'
$ws.Range('D23').Value = 0.8795596361160278
$ws.Range('E23').Value = 0.5047705173492432
$ws.Range('F23').Value = 0.5202772617340088
$ws.Range('G23').Value = 0.9378484078549304
$ws.Range('H23').Value = 25.02793967723846
$ws.Range('J23').Value = 'Random'
$ws.Range('K23').Value = '{''module__num_units'': 10, ''module__activation_func'': ReLU(), ''lr'': 0.1}'
$ws.Range('L23').Value = 0.8720538020133972
$ws.Range('M23').Value = 0.5051361918449402
$ws.Range('N23').Value = 0.5180037021636963
$ws.Range('O23').Value = 0.9338382097630173
$ws.Range('P23').Value = 24.99433308839798

# Hunk 13: sheet 'Option 1 - LR1 - DN2 (70-30)', row 23
$ws = $wb.Worksheets.Item('Option 1 - LR1 - DN2 (70-30)')
$ws.Range('B23').Value = 'Random'
$ws.Range('C23').Value = '{''module__num_units'': 20, ''module__activation_func'': Tanh(), ''lr'': 0.1}'
$ws.Range('D23').Value = 0.9536649584770203
$ws.Range('E23').Value = 0.5118158459663391
$ws.Range('F23').Value = 0.5077627897262573
$ws.Range('G23').Value = 0.9765577087284808
$ws.Range('H23').Value = 24.98691529035568
$ws.Range('J23').Value = 'Random'
$ws.Range('K23').Value = '{''module__num_units'': 10, ''module__activation_func'': ReLU(), ''lr'': 0.1}'
$ws.Range('L23').Value = 0.9326292276382446
$ws.Range('M23').Value = 0.5089412331581116
$ws.Range('N23').Value = 0.502464771270752
$ws.Range('O23').Value = 0.9657273050081191
$ws.Range('P23').Value = 24.97279196977615

# Hunk 14: sheet 'Option 1 - LR2 - DN1 (80-20)', row 22
$ws = $wb.Worksheets.Item('Option 1 - LR2 - DN1 (80-20)')
$ws.Range('B22').Value = 'Random'
$ws.Range('C22').Value = '{''module__num_units'': 20, ''module__activation_func'': Tanh(), ''lr'': 0.1}'
$ws.Range('D22').Value = 221.9242401123047
$ws.Range('E22').Value = 13.55620574951172
$ws.Range('F22').Value = 0.9227257370948792
$ws.Range('G22').Value = 14.89712187344605
$ws.Range('H22').Value = 170.1274991035461
$ws.Range('J22').Value = 'Random'
$ws.Range('K22').Value = '{''module__num_units'': 50, ''module__activation_func'': Tanh(), ''lr'': 0.1}'
$ws.Range('L22').Value = 221.9751586914062
$ws.Range('M22').Value = 13.55815315246582
$ws.Range('N22').Value = 0.9227643609046936
$ws.Range('O22').Value = 14.89883078269588
$ws.Range('P22').Value = 170.1757669448853

# Hunk 15: sheet 'Option 1 - LR2 - DN1 (60-40)', row 22
$ws = $wb.Worksheets.Item('Option 1 - LR2 - DN1 (60-40)')
$ws.Range('B22').Value = 'Random'
$ws.Range('C22').Value = '{''module__num_units'': 10, ''module__activation_func'': Tanh(), ''lr'': 0.02}'
$ws.Range('D22').Value = 239.0968475341797
$ws.Range('E22').Value = 15.46049022674561
$ws.Range('F22').Value = 69.87956237792969
$ws.Range('G22').Value = 15.46275678959543
$ws.Range('H22').Value = 185.5507731437683
$ws.Range('J22').Value = 'Random'
$ws.Range('K22').Value = '{''module__num_units'': 20, ''module__activation_func'': ReLU(), ''lr'': 0.02}'
$ws.Range('L22').Value = 167.5979766845703
$ws.Range('M22').Value = 12.94260406494141
$ws.Range('N22').Value = 58.66448211669922
$ws.Range('O22').Value = 12.94596372173854
$ws.Range('P22').Value = 182.9998135566711

# Hunk 16: sheet 'Option 1 - LR2 - DN1 (70-30)', row 22
$ws = $wb.Worksheets.Item('Option 1 - LR2 - DN1 (70-30)')
$ws.Range('B22').Value = 'Random'
$ws.Range('C22').Value = '{''module__num_units'': 10, ''module__activation_func'': Tanh(), ''lr'': 0.02}'
$ws.Range('D22').Value = 190.3064880371094
$ws.Range('E22').Value = 13.7652006149292
$ws.Range('F22').Value = 7.015673637390137
$ws.Range('G22').Value = 13.7951617619044
$ws.Range('H22').Value = 147.9785799980164
$ws.Range('J22').Value = 'Random'
$ws.Range('K22').Value = '{''module__num_units'': 10, ''module__activation_func'': ReLU(), ''lr'': 0.1}'
$ws.Range('L22').Value = 179.9839782714844
$ws.Range('M22').Value = 13.38486957550049
$ws.Range('N22').Value = 6.830113410949707
$ws.Range('O22').Value = 13.41581075714339
$ws.Range('P22').Value = 146.9086766242981
